$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# Fix typo in the "bgcolor" column (G): "#ccccccc" (7 c's) -> "#cccccc" (6 c's).
# Only the dropped-out / not-running / "Other" Republican rows (24-29, 31) get
# corrected; row 11 ("demother") keeps the original typo untouched, matching
# the target diff.
$rows = @(24, 25, 26, 27, 28, 29, 31)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "#cccccc"
}

# Update the active cell selection to H26
$ws.Range("H26").Select()
